$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("effort")

# Row 14: merge Effort + Additional Effort into single Effort value, drop Additional Effort
$ws.Range("B14").Value = 4
$ws.Range("C14").ClearContents()

# Row 37: merge Effort + Additional Effort into single Effort value, drop Additional Effort
$ws.Range("B37").Value = 2
$ws.Range("C37").ClearContents()

# New row 41: new log entry
$ws.Range("A40").Copy()
$ws.Range("A41").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A41").Value = 41227
$ws.Range("B41").Value = 2.5
$ws.Range("D41").Value = "Test case tc08 put to operation but still shows some bad behavior"

# Update the view: scroll back to top-left and move active selection to B22
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B22").Select()
